$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("Cariaini", $true, $false, $false, $false, $false, $true, 1, $false, "Caraiani", 2)
